# Simulated scheduled-task run: append the newest sensor reading (row 5)
# and refresh the precision of the previous row's timestamp (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh existing row 4 timestamp (tiny float precision update)
$ws.Range("A4").Value = 45869.50021532407

# New row 5 data
$ws.Range("A5").Value = 45869.62521966729
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 31
$ws.Range("D5").Value = 24.67
$ws.Range("E5").Value = 61.55
$ws.Range("F5").Value = 444.73
$ws.Range("G5").Value = 12.17
$ws.Range("H5").Value = "ESE"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "15:00:18"
